$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.432.51"
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").Value = "1.655.05"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.74"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.10"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").Value = "1.889.37"
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("D13").Value = "1.646.66"
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.87"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").Value = "27.437.79"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.31"
$ws.Range("E18").Value = "  -7.61%  "
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -3.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.32"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.62"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.94"
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.20"
$ws.Range("E30").Value = "  -4.12%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0497"
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("D33").Value = "1.463.62"
$ws.Range("E33").Value = "  +2.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("E35").Value = "  -4.30%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("E38").Value = "  -3.32%  "
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.51"
$ws.Range("E43").Value = "  -6.15%  "
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").Value = "1.797.79"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.782"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.73"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.39"
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  -5.27%  "
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("E51").Value = "  -0.17%  "
